# Apply 2021 MSRP updates to the product pricing sheet:
#  - Bump YEAR/BASE MSRP for the RC 300 / RC 350 lineup (rows 2-9) and RC F (row 53)
#  - Rename "RC F TRACK" -> "RC F FUJI SPEEDWAY EDITION" and update its year/price (row 54)
#  - Append four new "Black Line" special-edition trims as rows 95-98

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $year, $msrp) {
    $ws.Range("C$row").Value = $year
    $ws.Range("D$row").Value = $msrp
}

# RC 300 / RC 350 lineup: 2020 -> 2021, refreshed MSRPs
Set-Row 2 2021 42120
Set-Row 3 2021 46590
Set-Row 4 2021 44810
Set-Row 5 2021 48765
Set-Row 6 2021 45050
Set-Row 7 2021 49520
Set-Row 8 2021 47215
Set-Row 9 2021 51130

# RC F: 2020 -> 2021, refreshed MSRP
Set-Row 53 2021 65875

# New Black Line special-edition trims appended at the bottom of the table.
# Trim codes (column A) are entered first, then the model trim names
# (column B), then the year/price - matching the order the data was
# originally typed in.
$newRows = @(
    @{ Row = 95; Code = "9203SE"; Trim = "RC 300 F SPORT Black Line"; Year = 2021; Msrp = 49160 },
    @{ Row = 96; Code = "9207SE"; Trim = "RC 300 AWD F SPORT Black Line"; Year = 2021; Msrp = 51335 },
    @{ Row = 97; Code = "9213SE"; Trim = "RC 350 F SPORT Black Line"; Year = 2021; Msrp = 52090 },
    @{ Row = 98; Code = "9217SE"; Trim = "RC 350 AWD F SPORT Black Line"; Year = 2021; Msrp = 53700 }
)

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.Code
}
foreach ($r in $newRows) {
    $ws.Range("B$($r.Row)").Value = $r.Trim
}
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("C$row").Value = $r.Year
    $ws.Range("D$row").Value = $r.Msrp
    $ws.Range("D$row").NumberFormat = $ws.Range("D92").NumberFormat
    $ws.Range("E$row").Value = 1025
    $ws.Range("E$row").NumberFormat = $ws.Range("E92").NumberFormat
}

# RC F TRACK -> RC F FUJI SPEEDWAY EDITION, 2020 -> 2021, refreshed MSRP
$ws.Range("B54").Value = "RC F FUJI SPEEDWAY EDITION"
Set-Row 54 2021 97100

# Update the active selection to reflect where the editor left off
$ws.Range("D99").Select()
